$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2 "test1" -> "test2"; D2 "test2" -> boolean TRUE (WAHR/FALSCH format)
$ws.Range("C2").Value = "test2"
$ws.Range("D2").NumberFormat = '"WAHR";"WAHR";"FALSCH"'
$ws.Range("D2").Value = $true

# Row 3: B3 stays 200; D3 -> boolean FALSE
$ws.Range("D3").NumberFormat = '"WAHR";"WAHR";"FALSCH"'
$ws.Range("D3").Value = $false

# Row 4 (new): A4 "Test1", C4 "test3"
$ws.Range("A4").Value = "Test1"
$ws.Range("C4").Value = "test3"

# Row 5 (new): A5 "Test4", D5 -> boolean TRUE
$ws.Range("A5").Value = "Test4"
$ws.Range("D5").NumberFormat = '"WAHR";"WAHR";"FALSCH"'
$ws.Range("D5").Value = $true

# Row 6 (new): A6 "Test5", C6 "test6", D6 -> boolean FALSE
$ws.Range("A6").Value = "Test5"
$ws.Range("C6").Value = "test6"
$ws.Range("D6").NumberFormat = '"WAHR";"WAHR";"FALSCH"'
$ws.Range("D6").Value = $false

# Row 7 (new): B7 200, D7 -> boolean TRUE
$ws.Range("B7").Value = 200
$ws.Range("D7").NumberFormat = '"WAHR";"WAHR";"FALSCH"'
$ws.Range("D7").Value = $true

# Update selection to match target view state
[void]$ws.Range("F6").Select()
